# ----------------------------------------------------------------------------
# "add README.md, fix database" — rework the Sheet1 import layout:
#   - drop the NAMA KARYAWAN / TANGGAL* / NO PART / BARCODE / LPD / NAMA RAK
#     columns
#   - keep NOMOR POLISI / MODEL KENDARAAN / VIN RANGKA / KILOMETER (their
#     column formatting travels with them)
#   - insert two brand-new columns, MODEL and TIPE MESIN, with new sample
#     data
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (NAMA KARYAWAN) disappears entirely: deleting it shifts
# B..K (NOMOR POLISI..NAMA RAK) left by one, carrying their per-column
# width/bestFit metadata with them.
$ws.Columns("A").Delete()

# After that shift the columns we still want are A (NOMOR POLISI),
# B (MODEL KENDARAAN), C (VIN RANGKA) and D (KILOMETER); E..J
# (TANGGAL PERBAIKAN..NAMA RAK) are no longer needed.
$ws.Columns("E:J").Delete()

# Make room for the two new columns (MODEL, TIPE MESIN) between
# NOMOR POLISI and MODEL KENDARAAN.
$ws.Columns("B:C").Insert()

# ---- header row ----
$ws.Range("A1").Value = "NOMOR POLISI"
$ws.Range("B1").Value = "MODEL"
$ws.Range("C1").Value = "TIPE MESIN"
$ws.Range("D1").Value = "MODEL KENDARAAN"
$ws.Range("E1").Value = "VIN RANGKA"
$ws.Range("F1").Value = "KILOMETER"

# The KILOMETER value cell used to carry a date numeric format (style
# index 1 in the original file) - strip that back to General before
# writing the new plain number into it.
$ws.Range("F2").ClearFormats()

# ---- data row ----
$ws.Range("A2").Value = "D6242zah"
$ws.Range("B2").Value = "d1525"
$ws.Range("C2").Value = "yahaasd"
$ws.Range("D2").Value = "ahdaman"
$ws.Range("E2").Value = 1920
$ws.Range("F2").Value = 12466

# The two freshly-inserted columns (B, C) don't have a bestFit width yet -
# approximate Excel's "best fit" sizing for their new short headers
# (closest size this host can express is in 1/6-character increments).
$ws.Columns("B").ColumnWidth = 6.5
$ws.Columns("C").ColumnWidth = 10

# Reset the view: back to 100% zoom (was 85%) and the active cell
# lands on the last data cell instead of the old F5 selection.
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("F2").Select()
